$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$years = 2003..2023

$row = 3
foreach ($year in $years) {
    $path = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\$year.xlsx"
    $ws.Cells.Item($row, 1).Value = $path
    $ws.Cells.Item($row, 2).Value = "$year"
    $row = $row + 1
}

$ws.Range("A2:B23").Select()
